$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# NOTE: Column D ("Price") values are stored as plain TEXT in the source data
# (they use "." as a thousands separator, e.g. "35.572.92", which is not a valid
# number). Several of the new values (e.g. "1.01", "233.53") DO look like valid
# numbers to Excel, so a plain Range.Value assignment would silently convert them
# to numeric cells (losing formatting like trailing zeros). To faithfully reproduce
# a text update, we force the cell to Text format before writing, then restore the
# default "Normal" style so we do not leave a stray number-format behind.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Row 2
Set-TextValue $ws.Range("D2") '35.574.95'
$ws.Range("E2").Value = '  +2.38%  '

# Row 3
Set-TextValue $ws.Range("D3") '1.847.56'
$ws.Range("E3").Value = '  +1.98%  '

# Row 4
Set-TextValue $ws.Range("D4") '1.01'
$ws.Range("E4").Value = '  +0.83%  '

# Row 5
Set-TextValue $ws.Range("D5") '233.53'
$ws.Range("E5").Value = '  +3.79%  '

# Row 6
Set-TextValue $ws.Range("D6") '0.626'
$ws.Range("E6").Value = '  +3.20%  '

# Row 7
$ws.Range("E7").Value = '  +0.72%  '

# Row 8
Set-TextValue $ws.Range("D8") '44.17'
$ws.Range("E8").Value = '  +12.06%  '

# Row 9
Set-TextValue $ws.Range("D9") '0.313'
$ws.Range("E9").Value = '  +7.81%  '

# Row 10
$ws.Range("E10").Value = '  +3.97%  '

# Row 11
$ws.Range("E11").Value = '  +0.73%  '

# Row 12
Set-TextValue $ws.Range("D12") '2.114.92'
$ws.Range("E12").Value = '  +2.30%  '

# Row 13
Set-TextValue $ws.Range("D13") '11.45'
$ws.Range("E13").Value = '  +4.17%  '

# Row 14
$ws.Range("B14").Value = 'WrappedEther'
$ws.Range("C14").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextValue $ws.Range("D14") '1.859.22'
$ws.Range("E14").Value = '  +2.88%  '

# Row 15
$ws.Range("B15").Value = 'Polygon'
$ws.Range("C15").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
Set-TextValue $ws.Range("D15") '0.676'
$ws.Range("E15").Value = '  +6.48%  '

# Row 16
Set-TextValue $ws.Range("D16") '4.78'
$ws.Range("E16").Value = '  +8.70%  '

# Row 17
Set-TextValue $ws.Range("D17") '35.551.89'
$ws.Range("E17").Value = '  +2.58%  '

# Row 18
Set-TextValue $ws.Range("D18") '70.73'
$ws.Range("E18").Value = '  +3.54%  '

# Row 19
Set-TextValue $ws.Range("D19") '0.0₃0807'
$ws.Range("E19").Value = '  +4.90%  '

# Row 20
Set-TextValue $ws.Range("D20") '244.30'
$ws.Range("E20").Value = '  +1.06%  '

# Row 21
Set-TextValue $ws.Range("D21") '12.16'
$ws.Range("E21").Value = '  +9.10%  '

# Row 22
Set-TextValue $ws.Range("D22") '4.66'
$ws.Range("E22").Value = '  +13.36%  '

# Row 23
$ws.Range("E23").Value = '  +0.52%  '

# Row 24
Set-TextValue $ws.Range("D24") '2.25'
$ws.Range("E24").Value = '  +3.14%  '

# Row 25
Set-TextValue $ws.Range("D25") '171.68'
$ws.Range("E25").Value = '  +0.35%  '

# Row 26
Set-TextValue $ws.Range("D26") '7.98'
$ws.Range("E26").Value = '  +3.34%  '

# Row 27
Set-TextValue $ws.Range("D27") '17.90'
$ws.Range("E27").Value = '  +1.52%  '

# Row 28
$ws.Range("E28").Value = '  +0.94%  '

# Row 29
Set-TextValue $ws.Range("D29") '1.65'
$ws.Range("E29").Value = '  +34.61%  '

# Row 30
Set-TextValue $ws.Range("D30") '1.01'
$ws.Range("E30").Value = '  +0.79%  '

# Row 31
$ws.Range("B31").Value = 'Hedera'
$ws.Range("C31").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue $ws.Range("D31") '0.0568'
$ws.Range("E31").Value = '  +10.53%  '

# Row 32
$ws.Range("B32").Value = 'EURNeutrino'
$ws.Range("C32").Value = 'https://coinranking.com/coin/7YKHKSdb+eurneutrino-eurn'
Set-TextValue $ws.Range("D32") '3.342.24'
$ws.Range("E32").Value = '  +37.56%  '

# Row 33
Set-TextValue $ws.Range("D33") '4.14'
$ws.Range("E33").Value = '  +7.32%  '

# Row 34
Set-TextValue $ws.Range("D34") '3.98'
$ws.Range("E34").Value = '  +5.82%  '

# Row 35
$ws.Range("E35").Value = '  +1.23%  '

# Row 36
Set-TextValue $ws.Range("D36") '95.32'
$ws.Range("E36").Value = '  +15.36%  '

# Row 37
Set-TextValue $ws.Range("D37") '0.698'
$ws.Range("E37").Value = '  +8.01%  '

# Row 38
Set-TextValue $ws.Range("D38") '1.13'
$ws.Range("E38").Value = '  +6.53%  '

# Row 39
$ws.Range("B39").Value = 'VeChain'
$ws.Range("C39").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue $ws.Range("D39") '0.0198'
$ws.Range("E39").Value = '  +5.09%  '

# Row 40
$ws.Range("B40").Value = 'Maker'
$ws.Range("C40").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
Set-TextValue $ws.Range("D40") '1.347.33'
$ws.Range("E40").Value = '  +2.45%  '

# Row 41
$ws.Range("E41").Value = '  +7.63%  '

# Row 42
$ws.Range("B42").Value = 'InjectiveProtocol'
$ws.Range("C42").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextValue $ws.Range("D42") '15.30'
$ws.Range("E42").Value = '  +7.74%  '

# Row 43
$ws.Range("B43").Value = 'RenderToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue $ws.Range("D43") '2.41'
$ws.Range("E43").Value = '  +2.53%  '

# Row 44
$ws.Range("E44").Value = '  +3.02%  '

# Row 45
$ws.Range("E45").Value = '  +0.86%  '

# Row 46
Set-TextValue $ws.Range("D46") '2.82'
$ws.Range("E46").Value = '  +0.42%  '

# Row 47
Set-TextValue $ws.Range("D47") '6.31'
$ws.Range("E47").Value = '  +10.32%  '

# Row 48
Set-TextValue $ws.Range("D48") '0.0516'
$ws.Range("E48").Value = '  -0.26%  '

# Row 49
Set-TextValue $ws.Range("D49") '2.024.27'
$ws.Range("E49").Value = '  +3.08%  '

# Row 50
$ws.Range("E50").Value = '  +0.54%  '

# Row 51
Set-TextValue $ws.Range("D51") '102.75'
$ws.Range("E51").Value = '  +0.66%  '
